# Apply "Natmi following Dr Hou advice" update to Spon2-Itgb1 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Spon2"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.436226
$ws.Cells.Item(2, 8).Value = 1.308678
$ws.Cells.Item(2, 9).Value = 0.02612337021374315
$ws.Cells.Item(2, 10).Value = 0.02612337021374314
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 153.5290173333333
$ws.Cells.Item(2, 14).Value = 460.587052
$ws.Cells.Item(2, 15).Value = 0.3172206968818489
$ws.Cells.Item(2, 16).Value = 0.317220696881849
$ws.Cells.Item(2, 17).Value = 66.97334911525066
$ws.Cells.Item(2, 18).Value = 602.7601420372559
$ws.Cells.Item(2, 19).Value = 0.008286873704106135
$ws.Cells.Item(2, 20).Value = 0.008286873704106136

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Spon2"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.436226
$ws.Cells.Item(3, 8).Value = 1.308678
$ws.Cells.Item(3, 9).Value = 0.02612337021374315
$ws.Cells.Item(3, 10).Value = 0.02612337021374314
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 168.7997026666667
$ws.Cells.Item(3, 14).Value = 506.3991080000001
$ws.Cells.Item(3, 15).Value = 0.3487728915577651
$ws.Cells.Item(3, 16).Value = 0.3487728915577651
$ws.Cells.Item(3, 17).Value = 73.63481909546934
$ws.Cells.Item(3, 18).Value = 662.7133718592241
$ws.Cells.Item(3, 19).Value = 0.009111123366681189
$ws.Cells.Item(3, 20).Value = 0.009111123366681188

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Spon2"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.436226
$ws.Cells.Item(4, 8).Value = 1.308678
$ws.Cells.Item(4, 9).Value = 0.02612337021374315
$ws.Cells.Item(4, 10).Value = 0.02612337021374314
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 68.09032333333333
$ws.Cells.Item(4, 14).Value = 204.27097
$ws.Cells.Item(4, 15).Value = 0.1406878008722904
$ws.Cells.Item(4, 16).Value = 0.1406878008722904
$ws.Cells.Item(4, 17).Value = 29.70276938640666
$ws.Cells.Item(4, 18).Value = 267.32492447766
$ws.Cells.Item(4, 19).Value = 0.003675239506744218
$ws.Cells.Item(4, 20).Value = 0.003675239506744218

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Spon2"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.436226
$ws.Cells.Item(5, 8).Value = 1.308678
$ws.Cells.Item(5, 9).Value = 0.02612337021374315
$ws.Cells.Item(5, 10).Value = 0.02612337021374314
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 93.562673
$ws.Cells.Item(5, 14).Value = 280.688019
$ws.Cells.Item(5, 15).Value = 0.1933186106880956
$ws.Cells.Item(5, 16).Value = 0.1933186106880956
$ws.Cells.Item(5, 17).Value = 40.814470592098
$ws.Cells.Item(5, 18).Value = 367.330235328882
$ws.Cells.Item(5, 19).Value = 0.005050133636211605
$ws.Cells.Item(5, 20).Value = 0.005050133636211604

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Spon2"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 15.76143266666667
$ws.Cells.Item(6, 8).Value = 47.284298
$ws.Cells.Item(6, 9).Value = 0.9438725354525366
$ws.Cells.Item(6, 10).Value = 0.9438725354525365
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 153.5290173333333
$ws.Cells.Item(6, 14).Value = 460.587052
$ws.Cells.Item(6, 15).Value = 0.3172206968818489
$ws.Cells.Item(6, 16).Value = 0.317220696881849
$ws.Cells.Item(6, 17).Value = 2419.837269078832
$ws.Cells.Item(6, 18).Value = 21778.53542170949
$ws.Cells.Item(6, 19).Value = 0.2994159034638913
$ws.Cells.Item(6, 20).Value = 0.2994159034638913

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Spon2"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 15.76143266666667
$ws.Cells.Item(7, 8).Value = 47.284298
$ws.Cells.Item(7, 9).Value = 0.9438725354525366
$ws.Cells.Item(7, 10).Value = 0.9438725354525365
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 168.7997026666667
$ws.Cells.Item(7, 14).Value = 506.3991080000001
$ws.Cells.Item(7, 15).Value = 0.3487728915577651
$ws.Cells.Item(7, 16).Value = 0.3487728915577651
$ws.Cells.Item(7, 17).Value = 2660.525147734021
$ws.Cells.Item(7, 18).Value = 23944.72632960619
$ws.Cells.Item(7, 19).Value = 0.3291971534517403
$ws.Cells.Item(7, 20).Value = 0.3291971534517403

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Spon2"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 15.76143266666667
$ws.Cells.Item(8, 8).Value = 47.284298
$ws.Cells.Item(8, 9).Value = 0.9438725354525366
$ws.Cells.Item(8, 10).Value = 0.9438725354525365
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 68.09032333333333
$ws.Cells.Item(8, 14).Value = 204.27097
$ws.Cells.Item(8, 15).Value = 0.1406878008722904
$ws.Cells.Item(8, 16).Value = 0.1406878008722904
$ws.Cells.Item(8, 17).Value = 1073.201046469896
$ws.Cells.Item(8, 18).Value = 9658.809418229061
$ws.Cells.Item(8, 19).Value = 0.1327913513165703
$ws.Cells.Item(8, 20).Value = 0.1327913513165703

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Spon2"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 15.76143266666667
$ws.Cells.Item(9, 8).Value = 47.284298
$ws.Cells.Item(9, 9).Value = 0.9438725354525366
$ws.Cells.Item(9, 10).Value = 0.9438725354525365
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 93.562673
$ws.Cells.Item(9, 14).Value = 280.688019
$ws.Cells.Item(9, 15).Value = 0.1933186106880956
$ws.Cells.Item(9, 16).Value = 0.1933186106880956
$ws.Cells.Item(9, 17).Value = 1474.681770602851
$ws.Cells.Item(9, 18).Value = 13272.13593542566
$ws.Cells.Item(9, 19).Value = 0.1824681272203347
$ws.Cells.Item(9, 20).Value = 0.1824681272203346

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Spon2"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5010290000000001
$ws.Cells.Item(10, 8).Value = 1.503087
$ws.Cells.Item(10, 9).Value = 0.0300040943337204
$ws.Cells.Item(10, 10).Value = 0.0300040943337204
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 153.5290173333333
$ws.Cells.Item(10, 14).Value = 460.587052
$ws.Cells.Item(10, 15).Value = 0.3172206968818489
$ws.Cells.Item(10, 16).Value = 0.317220696881849
$ws.Cells.Item(10, 17).Value = 76.92249002550267
$ws.Cells.Item(10, 18).Value = 692.302410229524
$ws.Cells.Item(10, 19).Value = 0.009517919713851521
$ws.Cells.Item(10, 20).Value = 0.009517919713851519

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Spon2"
$ws.Cells.Item(11, 3).Value = "Itgb1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.5010290000000001
$ws.Cells.Item(11, 8).Value = 1.503087
$ws.Cells.Item(11, 9).Value = 0.0300040943337204
$ws.Cells.Item(11, 10).Value = 0.0300040943337204
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 168.7997026666667
$ws.Cells.Item(11, 14).Value = 506.3991080000001
$ws.Cells.Item(11, 15).Value = 0.3487728915577651
$ws.Cells.Item(11, 16).Value = 0.3487728915577651
$ws.Cells.Item(11, 17).Value = 84.57354622737735
$ws.Cells.Item(11, 18).Value = 761.1619160463962
$ws.Cells.Item(11, 19).Value = 0.01046461473934362
$ws.Cells.Item(11, 20).Value = 0.01046461473934362

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Spon2"
$ws.Cells.Item(12, 3).Value = "Itgb1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.5010290000000001
$ws.Cells.Item(12, 8).Value = 1.503087
$ws.Cells.Item(12, 9).Value = 0.0300040943337204
$ws.Cells.Item(12, 10).Value = 0.0300040943337204
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 68.09032333333333
$ws.Cells.Item(12, 14).Value = 204.27097
$ws.Cells.Item(12, 15).Value = 0.1406878008722904
$ws.Cells.Item(12, 16).Value = 0.1406878008722904
$ws.Cells.Item(12, 17).Value = 34.11522660937667
$ws.Cells.Item(12, 18).Value = 307.03703948439
$ws.Cells.Item(12, 19).Value = 0.004221210048975872
$ws.Cells.Item(12, 20).Value = 0.004221210048975872

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Spon2"
$ws.Cells.Item(13, 3).Value = "Itgb1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.5010290000000001
$ws.Cells.Item(13, 8).Value = 1.503087
$ws.Cells.Item(13, 9).Value = 0.0300040943337204
$ws.Cells.Item(13, 10).Value = 0.0300040943337204
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 93.562673
$ws.Cells.Item(13, 14).Value = 280.688019
$ws.Cells.Item(13, 15).Value = 0.1933186106880956
$ws.Cells.Item(13, 16).Value = 0.1933186106880956
$ws.Cells.Item(13, 17).Value = 46.87761249051701
$ws.Cells.Item(13, 18).Value = 421.898512414653
$ws.Cells.Item(13, 19).Value = 0.005800349831549391
$ws.Cells.Item(13, 20).Value = 0.00580034983154939

